$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The codeforiati:group-code and codeforiati:group-name columns (C and D)
# have swapped places: column C now holds the group name, and column D
# now holds the group code.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $colC = $ws.Cells.Item($r, 3).Value2
    $colD = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value = $colD
    $ws.Cells.Item($r, 4).Value = $colC
}
